# "Añadida ronda 6 datos de prueba"
# On the "Jueces" sheet, a new "Ronda 6" results column is inserted
# between the existing "Ronda 5" column (I) and "Final" column (K):
#   - J1 ("Semis") becomes "Ronda 6"
#   - a new column L is added, carrying what used to be in J ("Semis"),
#     with the same per-row "X" marks as column K (Final)
#   - K ("Final") is untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jueces")

# New column L: same per-row marks as column K (rows 1-25, i.e. header +
# the 24 judge rows that have data), including formatting/styles.
$ws.Range("K1:K25").Copy($ws.Range("L1")) | Out-Null

# Column L header keeps the "Semis" label that used to live in J1.
$ws.Range("L1").Value = "Semis"

# Column J header becomes the new "Ronda 6" label.
$ws.Range("J1").Value = "Ronda 6"

# Match the width of the new column to the column it inherited the
# "Semis" heading from.
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(10).ColumnWidth

# Final on-screen selection left by the edit.
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
